$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per the source diff.
# Numeric-looking text values are prefixed with a literal single-quote
# (PowerShell-escaped) to force Excel to keep them as text, matching the
# inlineStr representation in the original workbook, rather than silently
# converting them to numbers.

$ws.Range("D2").Value = '67.072.97'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '3.123.15'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '''174.98'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.122.86'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '''0.516'
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").Value = '''6.40'
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").Value = '''0.152'
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("E13").Value = '  -3.61%  '
$ws.Range("D14").Value = '''36.15'
$ws.Range("E14").Value = '  -2.97%  '
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").Value = '3.642.66'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '67.021.60'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").Value = '''17.09'
$ws.Range("E18").Value = '  +3.81%  '
$ws.Range("E19").Value = '  -1.64%  '
$ws.Range("D20").Value = '3.119.71'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '''485.79'
$ws.Range("E21").Value = '  -1.09%  '
$ws.Range("D22").Value = '''7.87'
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").Value = '''83.91'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").Value = '''12.81'
$ws.Range("E25").Value = '  -3.18%  '
$ws.Range("D27").Value = '''10.20'
$ws.Range("E27").Value = '  -1.56%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '''8.03'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("D33").Value = '''0.112'
$ws.Range("E33").Value = '  -1.93%  '
$ws.Range("D34").Value = '0.0₃0943'
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '''48.00'
$ws.Range("E36").Value = '  +1.61%  '
$ws.Range("D37").Value = '''5.62'
$ws.Range("E37").Value = '  -4.51%  '
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("D39").Value = '''0.312'
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").Value = '''49.28'
$ws.Range("E40").Value = '  -1.75%  '
$ws.Range("E41").Value = '  -0.43%  '
$ws.Range("E42").Value = '  -3.56%  '
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("D45").Value = '2.809.28'
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").Value = '''374.26'
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0349'
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("D48").Value = '''134.79'
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D50").Value = '''24.69'
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("E51").Value = '  +0.36%  '
